# Update odds values for the rows that changed (FlashScore weekly games file).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 updates
$ws.Range("G5").Value  = 2.45
$ws.Range("I5").Value  = 3.2
$ws.Range("J5").Value  = 3.25
$ws.Range("L5").Value  = 4
$ws.Range("M5").Value  = 1.1
$ws.Range("N5").Value  = 7
$ws.Range("Q5").Value  = 2.5
$ws.Range("R5").Value  = 1.5
$ws.Range("X5").Value  = 10
$ws.Range("Y5").Value  = 10
$ws.Range("AI5").Value = 15
$ws.Range("AW5").Value = 5

# Row 39 updates
$ws.Range("G39").Value  = 2.25
$ws.Range("I39").Value  = 3.1
$ws.Range("J39").Value  = 3
$ws.Range("K39").Value  = 2.1
$ws.Range("L39").Value  = 3.6
$ws.Range("M39").Value  = 1.06
$ws.Range("N39").Value  = 10
$ws.Range("O39").Value  = 1.3
$ws.Range("P39").Value  = 3.4
$ws.Range("Q39").Value  = 2.03
$ws.Range("R39").Value  = 1.83
$ws.Range("U39").Value  = 1.8
$ws.Range("V39").Value  = 1.91
$ws.Range("Y39").Value  = 9.5
$ws.Range("AA39").Value = 19
$ws.Range("AB39").Value = 29
$ws.Range("AF39").Value = 51
$ws.Range("AG39").Value = 251
$ws.Range("AH39").Value = 9.5
$ws.Range("AI39").Value = 15
$ws.Range("AJ39").Value = 11
$ws.Range("AL39").Value = 23
$ws.Range("AO39").Value = 13
$ws.Range("AP39").Value = 23
$ws.Range("AR39").Value = 67
$ws.Range("BB39").Value = 201
